$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 265, shifting existing rows 265-323 down to 266-324
$ws.Rows.Item(265).Insert()

# Fill in the new row 265 with the new data
$ws.Cells.Item(265, 1).Value = 10
$ws.Cells.Item(265, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(265, 3).Value = "La Araucanía"
$ws.Cells.Item(265, 4).Value = 44855
$ws.Cells.Item(265, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(265, 5).Value = 9
$ws.Cells.Item(265, 6).Value = "Fruta"
$ws.Cells.Item(265, 7).Value = 100102
$ws.Cells.Item(265, 8).Value = "Cítricos"
$ws.Cells.Item(265, 9).Value = 100102006
$ws.Cells.Item(265, 10).Value = "Pomelo"
$ws.Cells.Item(265, 11).Value = "Start Ruby"
$ws.Cells.Item(265, 12).Value = "Primera"
$ws.Cells.Item(265, 13).Value = 55
$ws.Cells.Item(265, 14).Value = 16000
$ws.Cells.Item(265, 15).Value = 16000
$ws.Cells.Item(265, 16).Value = 16000
$ws.Cells.Item(265, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(265, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(265, 19).Value = 1143
$ws.Cells.Item(265, 20).Value = 14
